$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 2-4 (2004, 2008, 2009 year rows), keeping row 1 (headers)
# and what was row 5 (2011 data) which will shift up to row 2.
$ws.Range("2:4").Delete()
